$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data for the new rows to append (Item, Multiplicador, Pontos, Batalhas, ValorApostado, ValorAcumulado, Resultado)
$newRows = @(
    @("BonusPower", 0.16, 0,  68,  40.3, 0,  "lose"),
    @("SkipBoss",   0,    0,  11,  20,   0,  "lose"),
    @("SkipBoss",   2,    10, 207, 40,   80, "win")
)

$startRow = 40

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
